$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto "Price" column (D) holds values formatted as plain-text strings
# (e.g. with thousand-separator dots, or trailing zeros that are significant).
# Excel auto-converts plain numeric-looking input into a real number, which
# would silently strip meaningful trailing zeros (e.g. "18.10" -> 18.1).
# Force those specific target cells to Text format first so the literal
# string is preserved exactly, matching how the source data is stored.
$textCells = @("D5", "D7", "D9", "D10", "D12", "D14", "D15", "D16", "D19", "D20", "D22", "D24", "D25", "D26", "D27", "D30", "D31", "D36", "D38", "D39", "D40", "D41", "D45", "D46", "D47", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '38.772.33'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '2.103.30'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '227.96'
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("D7").Value = '62.37'
$ws.Range("E7").Value = '  +1.37%  '
$ws.Range("D9").Value = '0.390'
$ws.Range("E9").Value = '  +1.78%  '
$ws.Range("D10").Value = '0.0841'
$ws.Range("E10").Value = '  -0.66%  '
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("D12").Value = '15.74'
$ws.Range("E12").Value = '  +6.00%  '
$ws.Range("D13").Value = '2.415.72'
$ws.Range("E13").Value = '  +0.25%  '
$ws.Range("D14").Value = '22.04'
$ws.Range("E14").Value = '  -1.72%  '
$ws.Range("D15").Value = '0.808'
$ws.Range("E15").Value = '  +2.95%  '
$ws.Range("D16").Value = '5.52'
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("D17").Value = '2.110.06'
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").Value = '38.767.36'
$ws.Range("E18").Value = '  +0.44%  '
$ws.Range("D19").Value = '71.93'
$ws.Range("E19").Value = '  +1.20%  '
$ws.Range("D20").Value = '6.11'
$ws.Range("E20").Value = '  +0.70%  '
$ws.Range("D21").Value = '0.0₃0840'
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").Value = '227.63'
$ws.Range("E22").Value = '  +0.54%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").Value = '2.34'
$ws.Range("E24").Value = '  -2.69%  '
$ws.Range("D25").Value = '2.33'
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("D26").Value = '9.64'
$ws.Range("E26").Value = '  +1.89%  '
$ws.Range("D27").Value = '172.44'
$ws.Range("E27").Value = '  +0.93%  '
$ws.Range("E28").Value = '  +4.46%  '
$ws.Range("E29").Value = '  +4.72%  '
$ws.Range("D30").Value = '19.35'
$ws.Range("E30").Value = '  +1.15%  '
$ws.Range("D31").Value = '2.51'
$ws.Range("E31").Value = '  +9.97%  '
$ws.Range("E32").Value = '  +0.28%  '
$ws.Range("E33").Value = '  +1.36%  '
$ws.Range("E34").Value = '  -0.33%  '
$ws.Range("E35").Value = '  +7.06%  '
$ws.Range("D36").Value = '0.0620'
$ws.Range("E36").Value = '  +2.07%  '
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("D38").Value = '3.58'
$ws.Range("E38").Value = '  +0.25%  '
$ws.Range("D39").Value = '0.998'
$ws.Range("E39").Value = '  -0.20%  '
$ws.Range("D40").Value = '18.10'
$ws.Range("E40").Value = '  -3.14%  '
$ws.Range("D41").Value = '102.73'
$ws.Range("E41").Value = '  +2.69%  '
$ws.Range("E42").Value = '  +3.48%  '
$ws.Range("D43").Value = '1.526.54'
$ws.Range("E43").Value = '  -1.25%  '
$ws.Range("E44").Value = '  +6.67%  '
$ws.Range("D45").Value = '2.81'
$ws.Range("E45").Value = '  -0.90%  '
$ws.Range("D46").Value = '7.80'
$ws.Range("E46").Value = '  +1.61%  '
$ws.Range("D47").Value = '0.0911'
$ws.Range("E47").Value = '  -0.40%  '
$ws.Range("D48").Value = '4.14'
$ws.Range("E48").Value = '  -0.97%  '
$ws.Range("E49").Value = '  +1.63%  '
$ws.Range("E50").Value = '  -0.47%  '
$ws.Range("D51").Value = '2.301.81'
$ws.Range("E51").Value = '  +0.27%  '
